# Insert a new row at position 167, shifting existing rows 167-172 down to 168-173.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(167).Insert()

# Populate the newly inserted row 167 with the new weekly record.
$ws.Range("A167").Value2 = 4
$ws.Range("B167").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C167").Value2 = "Los Lagos"
$ws.Range("D167").Value2 = 44509
$ws.Range("E167").Value2 = 10
$ws.Range("F167").Value2 = 100112043
$ws.Range("G167").Value2 = "Pepino ensalada"
$ws.Range("H167").Value2 = "Sin especificar"
$ws.Range("I167").Value2 = "Primera"
$ws.Range("J167").Value2 = 400
$ws.Range("K167").Value2 = 10000
$ws.Range("L167").Value2 = 12000
$ws.Range("M167").Value2 = 11000
$ws.Range("N167").Value2 = "`$/caja 60 unidades"
$ws.Range("O167").Value2 = "Región de Arica y Parinacota"
$ws.Range("P167").Value2 = 183
$ws.Range("Q167").Value2 = 60
$ws.Range("R167").Value2 = "Hortaliza"
